$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values for B/C/D columns, rows 3-17 (Data_carac_2, Data_carac_3, Data_carac_4)
$data = @{
    3  = @(8932.4666666666672, 0.81458949777723322, 9051.5333333333328)
    4  = @(6715.9333333333334, 0.34018086301823247, 9250.4666666666672)
    5  = @(7165.1999999999998, 0.3041257114962968, 11035.4)
    6  = @(7989.5333333333338, 0.17550673043891526, 18996.333333333332)
    7  = @(13822.799999999999, 0.42419100321698655, 19708.933333333334)
    8  = @(17924.799999999999, 0.64695149929858065, 20164.599999999999)
    9  = @(10624.666666666666, 0.51569465887054156, 13374.533333333333)
    10 = @(14527.933333333332, 0.55620997233010594, 21383.533333333333)
    11 = @(22070.200000000001, 0.73979512088145361, 24145.599999999999)
    12 = @(16753.133333333335, 0.71573005781973198, 19591.733333333334)
    13 = @(16583.866666666665, 0.40567172762762926, 27256.133333333335)
    14 = @(12498.266666666666, 0.39298358244176668, 18932.400000000001)
    15 = @(17679.200000000001, 0.45222541611648842, 21893.866666666665)
    16 = @(24116.066666666666, 0.48426089877576123, 31433.200000000001)
    17 = @(20331.466666666667, 0.43050183458027658, 30325.799999999999)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}
